# "add thread executor base"
# Adds a few new "open source" library candidates (liblfds, theron, mintomic)
# to the "open source" sheet, and makes that sheet the active/selected tab.

$wb = $excel.ActiveWorkbook
$wsOpenSource = $wb.Worksheets.Item("open source")

# New row 7: liblfds-6.1.1 candidate entry (left block: name / license / description)
$wsOpenSource.Range("A7").Value2 = "liblfds-6.1.1"
$wsOpenSource.Range("B7").Value2 = "FREE"
$wsOpenSource.Range("C7").Value2 = "免锁数据结构"

# New row 15: theron candidate entry (right block: alternative name / description)
$wsOpenSource.Range("E15").Value2 = "theron"
$wsOpenSource.Range("G15").Value2 = "并发操作"

# New row 16: mintomic candidate entry (right block: alternative name / license / description)
$wsOpenSource.Range("E16").Value2 = "mintomic"
$wsOpenSource.Range("F16").Value2 = "？？？"
$wsOpenSource.Range("G16").Value2 = "免锁数据结构"

# Move the selection/active cell & make "open source" the active tab
$wsOpenSource.Range("B18").Select()
$wsOpenSource.Activate()
